# Apply the "Push And Commit as on 04-May-2019" edit:
#  1. Remove the extra "PatientNotBrought"/"yes" column (old column F) from the
#     "PanelApp" sheet, shifting "IsCorporateAppointment"/"no" (old column G)
#     left into column F.
#  2. Append a new, empty worksheet named "Sheet3" at the end of the workbook.
#  3. Update the saved cursor/selection on "DoctorApp" (F10) and "PanelApp" (I7),
#     keeping "PanelApp" as the active/selected tab, as it was originally.

$wb = $excel.ActiveWorkbook

$doctorApp = $wb.Worksheets.Item("DoctorApp")
$panelApp  = $wb.Worksheets.Item("PanelApp")

# --- 1. Delete the obsolete "PatientNotBrought" column on PanelApp ----------
# This removes old column F entirely; old column G ("IsCorporateAppointment")
# shifts left to become the new column F.
$panelApp.Columns.Item(6).Delete()

# --- 2. Add the new empty "Sheet3" worksheet at the very end ---------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet3"

# --- 3. Restore selections / active sheet -----------------------------------
[void]$doctorApp.Range("F10").Select()

[void]$panelApp.Activate()
[void]$panelApp.Range("I7").Select()
